$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column O: new "clean_description_2" column (header + 63 data rows)
$values = @(
    'clean_description_2',
    'amber color aroma peach jam orange blossom flavor mélange mellow fruit savory note particular apricot lemon curd quince bitter almond despite strong tannin moment brightness across entire palate linger finish',
    'soft fruity juicy pleasant simple blend different european country balance drinkable party barbecue screwcap',
    'violetred color semisweet nose blackberry vanilla soft entry flavor black cherry blueberry creme brûlée smooth palate blueberry pie finish brighten brief unwelcome tart note',
    'mouthwatering blend start nose fresh mint pineapple mango follow flavor orange peel peach lemon fresh fruity mineral backbone finish earthy spin pair grill sardine poultry dish',
    'orangestyle cloudy yellowgold appearance spend day skin day lee aroma flavor age apple cider tart apple pie astringent finish style everyone appreciate orange make chardonnay solid example',
    'blend melnik grenache noir petit verdot aroma saddle leather cassis dark chocolate mouth flavor cherry chocolate dry blueberry good balance soft tannic finish',
    'light fruity sweetness gentle series white fruit acidity lemon zest end screwcap',
    'furmint grow marl soil aroma brioche apple tart caramelize pear flavor lemonlime juice green apple lime blossom crisp clean finish',
    'jumpy jammy aroma foxy black fruit slightly volatile narrow mouthfeel thats high acidity blacken berry flavor spicy peppery tannat toasty livewire acidic finish',
    'winemaker bartho eksteen wooded sauvy show incredible restraint harmony balance open attractive note ripe gooseberry apple flesh green melon rind palate offer justripe tropicalfruit tone pleasantly pithy texture limerind flavor hang linger finish enjoy bright acidity though mature well',
    'delicate white flower spin lemon peel nose start refine white bright fruit palate tropical exotic minerality give lift fullerbodied poise age potential foodfriendly character',
    'aroma black berry dry red raspberry blueberry mouth flavor blueberry juice black plum finish slightly jammy yet persistent',
    'aroma green apple white flower prepare palate fill mouth soft fruit flavor particular apple grapefruit pineapple soft finish mark note granny smith apple',
    'balance aroma green herb citrus zest fresh vital aggressive lemony twinge acidity lend snap palate taste pithy briny tangerine citrus fruit flavor salt lime grapefruit share finish',
    'reasonably rich concentrate example pinot gris one seem bite repress aromatically hint rather suggest fig pear melon end weighty note',
    'interest blend indigenous bulgarian uncommon french variety red make equal part melnik mourvèdre marselan offer aroma red cherry cassis leather flavor suggest black cherry tart cherry pie fine tannin persistent finish',
    'subdue citrus pear note nose find fuller expression round textured palate almost seem pleasingly bitter edge earth bosc pear peel generous mediumbodied',
    'scent clover stem green herb red currant lead complex add sassafras highlight glide gentle finish ferment wild yeast age month onequarter new french oak wellmade though subtle expression variety',
    'crisp apple freshness almost tip full citrus zestiness nose subtle graceful slight spritz tickle palate textured dry palate lemony crisp light full',
    'blend gamay prokupe aroma tart cherry mediterranean herb flavor recall cherry black plum stew prune',
    'aroma honeysuckle lemon curd flavor lemonmeringue pie lemon pith finish clean crisp',
    'whiff citrus show restrain nose palate crisp even lean disclose slight spritz riper fruit wait core savor midpalate find light refresh',
    'basic cherry aroma turn earthy soupy air crisp pull palate hold flavor tomato spice plum mild nothing much speak finish',
    'dark color rich jammy bake aroma char oak lush berry open blend cabernet sauvignon carmenère cabernet franc moderately bold acidity firm tannin lend energy structure palate taste blackberry cassis spicy oak slightly hot finish lot wood spice close act',
    'blend comprise merlot cabernet sauvignon melnik petit verdot aroma powder cocoa black cherry flavor tart cherry blackberry blueberry pleasant burst acidity show finish',
    'entice aroma blueberry syrup open blend merlot cabernet sauvignon melnik syrah petit verdot flavor cherry tart cherry pie blueberry red raspberry come burst acidity persistent finish',
    'darkgarnet aroma eucalyptus balsa wood black cherry flavor cassis leather black cherry lead long finish',
    'aroma cherry blueberry rise petal prepare taste bud flavor cherry pie vanilla custard black pepper turkish delight tannin feel momentarily stiff soon outweigh bold acidity carry bright black cherry finish',
    'almost caramel color offer aroma lemon zest vanilla bean grippy tannin provide firm backdrop flavor lemon clementine green apple caramel may seem chalky drinker enjoy alongside nut hard cheese put balance focus',
    'winemaker gordon newton johnson pretty elegant pinot certainly would expect new world region express masterful balance south african pinot achieve new world style ripe redfruit character old world finesse bright acidity earthy forestfloor dominance press floral tone partner note red cherry cola spice sou bois throughout silky smooth offer soft fine tannin somehow still manage provide present structure overall balance harmony now–',
    'heady florality damask rise join powerful earthy force acquire rich autumnal touch ripe fruit fragrance palate wonderfully dry precise fresh enough acidity brighten every last corner ethereal earthy elegant firm fragrant fresh—from first whiff last vestige tongue',
    'make grill meat tobacco spice aroma combination fruit savory flavor assertive character finish spin cedar',
    'amber eye feature entice bouquet almond blossom toffee slightly chalky palate sophisticate mélange bake apple peach smoke meat pencillead flavor touch brightness mark soft finish',
    'bulgarian chardonnay ferment eggshaped concrete vessel aroma anise frond citrus blossom lively flavor lemon zest grapefruit pith lemon curd mouth nice balance pleasant touch ripe yellow peach finish',
    'dry furmint hungary balaton region aroma lemon blossom green apple mouth flavor crisp apple cider fresh pear',
    'vranac garnet color aroma black cherry black plum convey flavor fresh cherry juice freshly slice black plum touch pleasant astringency finish',
    'interest blend indigenous bulgarian uncommon french variety red make equal part melnik mourvèdre marselan offer aroma red cherry cassis leather flavor suggest black cherry tart cherry pie fine tannin persistent finish',
    'white blend comprise chardonnay gewürztraminer delightful aroma white flower honeysuckle lime pith flavor jasmine lemongrass dry rosemary mouth finish clean crisp refresh',
    'make grape grow lake balaton region aroma lemon blossom vanilla lemon curd lead way wellcombined flavor green apple pear pineapple medium body palate clean finish',
    'there pleasant fruity richness pronounce aroma flavor peach pear guava hit note toast bake spice generous palate full feel welcome acidic lift close',
    'make melnick garnetcolored aroma black cherry blueberry roll tobacco leaf mouth flavor cherry cassis blueberry firm tannic finish',
    'neutral nose lead onto applefresh palate golden delicious pearmain apple creamy aspect bath softer friendly light crisp round apple fruit charm throughout balance morillon local term chardonnay',
    'aroma cut grass lemongrass honeysuckle pave way flavor lemon curd pineapple juice freshly cut guava balance mouth pleasantly acidic finish',
    'mediumstrawcolored aroma citrus blossom lemon zest whiff sea mist open flavor peach apricot touch salt tip tongue get great minerality crisp finish',
    'aromatic hint tomato leaf pea shoot make highly appetizing nose concentrate body add lovely notion ripe mirabelle plum combine fulsome fruit ample zest midpalate maximize flavor leave refresh finish',
    'streak blossom citrus pear discern nose palate aroma take zesty bright pure citrus note slight spritz add notion utter briskness fresh easy clean',
    'organically produce comprise melnik grenache noir petit verdot aroma blueberry tobacco leaf leather flavor chocolate powder cassis blackberry long finish pleasant burst acidity end',
    'plum black cherry anise lead nose rich red vanilla mint cherry somewhat medicinal edge overall balance offer exotic edge',
    'pinot gris switzerland aroma white flower green apple lemon pith mouth flavor lemon curd guava medium body palate clean finish',
    'lively whiff white blossom muscat grape lemon enticingly streak nose flavor continue featherlight dry palate fresh acidity make shine even brighter express delicious dry totally aromatic lightness refresh make summer irresistible lemon finish',
    'aroma blackberry strawberry chocolate open flavor blackberry red delicious apple peel aniseed vanilla harmonious palate understate tannin breakthrough note summer peach finish',
    'reasonably rich concentrate example pinot gris one seem bite repress aromatically hint rather suggest fig pear melon end weighty note',
    'aroma crush stone graphite lead toasty char berry scent tannat feel fresh lively full body benefit brace acidity lightly herbal plum berry flavor bordeauxlike feel dry taste bite salty finish',
    'lively whiff white blossom muscat grape lemon enticingly streak nose flavor continue featherlight dry palate fresh acidity make shine even brighter express delicious dry totally aromatic lightness refresh make summer irresistible lemon finish',
    'amber eye feature entice bouquet almond blossom toffee slightly chalky palate sophisticate mélange bake apple peach smoke meat pencillead flavor touch brightness mark soft finish',
    'deep redviolet nose blackberry cassis roast meat strong tannin buffer equally strong acidity flavor cassis black cherry flame broil meat leather violet rigid tannin subdue bright burst cranberry finish',
    'dark ruby glass offer bouquet black cherry pomegranate paprika web soft tannin shore flavor cherry cranberry roast red pepper marjoram constant bright fruit palate mark touch fresh dry bell pepper culminate citrusdrenched finish',
    'cigar box cedar red fruit aroma start bale red famed winemaker yiannis voyatzis palate spicy profile lead good acidity structure base luscious dark cherry flavor serious lively spin',
    'swiss red blend compose four variety aroma black cherry red pepper bell pepper aroma black raspberry elderberry tannin velvety well balance',
    'dry spicy aroma dusty plum tomato add good whole cab crisp alert palate pop acidity create sharpness earthy spicy currant redplum flavor finish direct juicy',
    'el capricho one uruguay consistent winery term put character quality tempranillo smoky spicy woodsy nose note cola accent berry currant scent pull draw palate bite tannic taste toasty oak plum tomato',
    'blend syrah cabernet sauvignon merlot inky garnetcolored offer aroma cassis elderberry palate combination cook fruit cool spice flavor yet lack accompany acidity flavor cassis elderberry anise orange peel vanilla back silky tannin stay smooth cool finish',
    'offer delightful bouquet blackberry rise petal touch salinity smooth silky first sip flavor black plum cherry chocolate lavender freshly roast coffee make way around mouth stern tannin make presence know recede finish notable linger espresso bean'
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 15).Value = $values[$i]
}

# Match header formatting used by the rest of row 1 (bold, bordered, centered)
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("O1").Value = $values[0]

$wb.Save()
